$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "B" / "D" column data -------------------------------------------------
# Written in this specific order so that newly-introduced shared strings are
# appended to sharedStrings.xml in the same order as the target workbook.
$ws.Range("B11").Value = "BANGALORE"
$ws.Range("B9").Value  = "NOIDA"
$ws.Range("B1").Value  = "CHENNAI"
$ws.Range("B2").Value  = "ANY"
$ws.Range("B3").Value  = "HYDERABAD"
$ws.Range("B4").Value  = "HYDERABAD"
$ws.Range("B8").Value  = "BANGALORE- conditional"
$ws.Range("B10").Value = "conditional - bangalore"
$ws.Range("B6").Value  = "Delhi"
$ws.Range("B7").Value  = "Delhi"

$ws.Range("D1").Value  = "T"
$ws.Range("D2").Value  = "T"
$ws.Range("D6").Value  = "T"
$ws.Range("D9").Value  = "T"
$ws.Range("D11").Value = "T"

# --- Column A text change -------------------------------------------------------
$ws.Range("A5").Value = "Course5  intelligence customer segmentation"

# --- Row heights ------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(11).RowHeight = 45

# --- Column widths ------------------------------------------------------------
# (inputs chosen so the engine's char-width rounding lands on the target
#  stored widths of 34 and ~23.43)
$ws.Columns.Item(1).ColumnWidth = 33.14
$ws.Columns.Item(2).ColumnWidth = 22.65

# --- Selection / view -------------------------------------------------------
$ws.Range("B6").Select()
